# "menghapus perbedaan harga toko dan pelanggan"
# Collapse the separate "Harga Pelanggan" (G) and "Harga Toko" (H) columns
# into a single "Harga Jual" column, shifting the remaining columns
# (Agen, Keterangan) left by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Harga Pelanggan" column entirely; this shifts "Harga Toko"
# into column G, "Agen" into H, and "Keterangan" into I.
$ws.Columns("G").EntireColumn.Delete()

# Copy the header formatting from the neighbouring "Harga Modal" header
# onto the now-merged price column, then rename it to "Harga Jual".
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Harga Jual"

# Move the active selection to match the new layout.
[void]$ws.Range("G2").Select()
